# feat: add redirect homepage
# Replace the placeholder "1" sample-data row (row 2) with real data,
# then update the sheet view (zoom + active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Jesús Lara"
$ws.Range("E2").Value = "Valencia"
$ws.Range("G2").Value = "Instituto Pedrito"
$ws.Range("H2").Value = "Jose Ramirez"
$ws.Range("J2").Value = "a"

$excel.ActiveWindow.Zoom = 85
$ws.Range("E14").Select()
